$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.555934
$ws.Range("H2").Value = 1.667802
$ws.Range("I2").Value = 0.005745252779589096
$ws.Range("J2").Value = 0.005745252779589094
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2288403333333333
$ws.Range("N2").Value = 0.686521
$ws.Range("O2").Value = 0.001525084821945529
$ws.Range("P2").Value = 0.001525084821945529
$ws.Range("Q2").Value = 0.1272201218713334
$ws.Range("R2").Value = 1.144981096842
$ws.Range("S2").Value = 0.00000876199781239169
$ws.Range("T2").Value = 0.000008761997812391689
$ws.Range("G3").Value = 0.555934
$ws.Range("H3").Value = 1.667802
$ws.Range("I3").Value = 0.005745252779589096
$ws.Range("J3").Value = 0.005745252779589094
$ws.Range("M3").Value = 0.2676766666666667
$ws.Range("N3").Value = 0.80303
$ws.Range("O3").Value = 0.001783905903194393
$ws.Range("P3").Value = 0.001783905903194393
$ws.Range("Q3").Value = 0.1488105600066667
$ws.Range("R3").Value = 1.33929504006
$ws.Range("S3").Value = 0.00001024899034885298
$ws.Range("T3").Value = 0.00001024899034885298
$ws.Range("G4").Value = 0.555934
$ws.Range("H4").Value = 1.667802
$ws.Range("I4").Value = 0.005745252779589096
$ws.Range("J4").Value = 0.005745252779589094
$ws.Range("M4").Value = 149.554372
$ws.Range("N4").Value = 448.6631160000001
$ws.Range("O4").Value = 0.99669100927486
$ws.Range("P4").Value = 0.9966910092748601
$ws.Range("Q4").Value = 83.14236024344802
$ws.Range("R4").Value = 748.2812421910321
$ws.Range("S4").Value = 0.00572624179142785
$ws.Range("T4").Value = 0.005726241791427849
$ws.Range("I5").Value = 0.823525905561055
$ws.Range("J5").Value = 0.823525905561055
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.2288403333333333
$ws.Range("N5").Value = 0.686521
$ws.Range("O5").Value = 0.001525084821945529
$ws.Range("P5").Value = 0.001525084821945529
$ws.Range("Q5").Value = 18.23576265293078
$ws.Range("R5").Value = 164.121863876377
$ws.Range("S5").Value = 0.001255946859050112
$ws.Range("T5").Value = 0.001255946859050112
$ws.Range("I6").Value = 0.823525905561055
$ws.Range("J6").Value = 0.823525905561055
$ws.Range("M6").Value = 0.2676766666666667
$ws.Range("N6").Value = 0.80303
$ws.Range("O6").Value = 0.001783905903194393
$ws.Range("P6").Value = 0.001783905903194393
$ws.Range("Q6").Value = 21.33054121167889
$ws.Range("R6").Value = 191.97487090511
$ws.Range("S6").Value = 0.001469092724363874
$ws.Range("T6").Value = 0.001469092724363874
$ws.Range("I7").Value = 0.823525905561055
$ws.Range("J7").Value = 0.823525905561055
$ws.Range("M7").Value = 149.554372
$ws.Range("N7").Value = 448.6631160000001
$ws.Range("O7").Value = 0.99669100927486
$ws.Range("P7").Value = 0.9966910092748601
$ws.Range("Q7").Value = 11917.64577412832
$ws.Range("R7").Value = 107258.8119671549
$ws.Range("S7").Value = 0.8208008659776409
$ws.Range("T7").Value = 0.820800865977641
$ws.Range("G8").Value = 0.3446996666666666
$ws.Range("H8").Value = 1.034099
$ws.Range("I8").Value = 0.003562269474506148
$ws.Range("J8").Value = 0.003562269474506148
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.2288403333333333
$ws.Range("N8").Value = 0.686521
$ws.Range("O8").Value = 0.001525084821945529
$ws.Range("P8").Value = 0.001525084821945529
$ws.Range("Q8").Value = 0.07888118661988888
$ws.Range("R8").Value = 0.7099306795789999
$ws.Range("S8").Value = 0.000005432763107249201
$ws.Range("T8").Value = 0.000005432763107249201
$ws.Range("G9").Value = 0.3446996666666666
$ws.Range("H9").Value = 1.034099
$ws.Range("I9").Value = 0.003562269474506148
$ws.Range("J9").Value = 0.003562269474506148
$ws.Range("M9").Value = 0.2676766666666667
$ws.Range("N9").Value = 0.80303
$ws.Range("O9").Value = 0.001783905903194393
$ws.Range("P9").Value = 0.001783905903194393
$ws.Range("Q9").Value = 0.09226805777444444
$ws.Range("R9").Value = 0.8304125199699999
$ws.Range("S9").Value = 0.000006354753544340706
$ws.Range("T9").Value = 0.000006354753544340706
$ws.Range("G10").Value = 0.3446996666666666
$ws.Range("H10").Value = 1.034099
$ws.Range("I10").Value = 0.003562269474506148
$ws.Range("J10").Value = 0.003562269474506148
$ws.Range("M10").Value = 149.554372
$ws.Range("N10").Value = 448.6631160000001
$ws.Range("O10").Value = 0.99669100927486
$ws.Range("P10").Value = 0.9966910092748601
$ws.Range("Q10").Value = 51.55134217694267
$ws.Range("R10").Value = 463.962079592484
$ws.Range("S10").Value = 0.003550481957854558
$ws.Range("T10").Value = 0.003550481957854558
$ws.Range("G11").Value = 16.17571666666667
$ws.Range("H11").Value = 48.52715
$ws.Range("I11").Value = 0.1671665721848498
$ws.Range("J11").Value = 0.1671665721848498
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.2288403333333333
$ws.Range("N11").Value = 0.686521
$ws.Range("O11").Value = 0.001525084821945529
$ws.Range("P11").Value = 0.001525084821945529
$ws.Range("Q11").Value = 3.701656393905556
$ws.Range("R11").Value = 33.31490754515
$ws.Range("S11").Value = 0.0002549432019757761
$ws.Range("T11").Value = 0.0002549432019757761
$ws.Range("G12").Value = 16.17571666666667
$ws.Range("H12").Value = 48.52715
$ws.Range("I12").Value = 0.1671665721848498
$ws.Range("J12").Value = 0.1671665721848498
$ws.Range("M12").Value = 0.2676766666666667
$ws.Range("N12").Value = 0.80303
$ws.Range("O12").Value = 0.001783905903194393
$ws.Range("P12").Value = 0.001783905903194393
$ws.Range("Q12").Value = 4.329861918277778
$ws.Range("R12").Value = 38.9687572645
$ws.Range("S12").Value = 0.0002982094349373252
$ws.Range("T12").Value = 0.0002982094349373252
$ws.Range("G13").Value = 16.17571666666667
$ws.Range("H13").Value = 48.52715
$ws.Range("I13").Value = 0.1671665721848498
$ws.Range("J13").Value = 0.1671665721848498
$ws.Range("M13").Value = 149.554372
$ws.Range("N13").Value = 448.6631160000001
$ws.Range("O13").Value = 0.99669100927486
$ws.Range("P13").Value = 0.9966910092748601
$ws.Range("Q13").Value = 2419.149147733267
$ws.Range("R13").Value = 21772.3423295994
$ws.Range("S13").Value = 0.1666134195479367
$ws.Range("T13").Value = 0.1666134195479367
